$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Hoja1" -> "Inventario"
$ws.Name = "Inventario"

# Update unit labels used in column C: "und" -> "unidades", "ml" -> "gramos"
$usedRange = $ws.UsedRange
foreach ($cell in $usedRange.Cells) {
    $val = $cell.Value2
    if ($val -eq "und") {
        $cell.Value2 = "unidades"
    } elseif ($val -eq "ml") {
        $cell.Value2 = "gramos"
    }
}

# Move the active selection to D17
$ws.Range("D17").Select()
